$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B69: the politeness_score was stored as text "3"; fix it to a genuine number 3
$ws.Range("B69").Value = 3

# Append new annotation row 70
$ws.Range("A70").Value = "Ruilin"
# Keep B70 as text "4" (leading apostrophe forces text-storage for the
# numeric-looking string), then strip the resulting quote-prefix format
# so no extra cell style is introduced.
$ws.Range("B70").Value = "'4"
$ws.Range("B70").ClearFormats()
$ws.Range("C70").Value = "could be a game changer."
$ws.Range("D70").Value = "APC"
$ws.Range("E70").Value = "OTH"
$ws.Range("F70").Value = "a3c87a5d-b7d4-4eb1-9136-458357f6153b"
$ws.Range("G70").Value = "IrVvIL2BaXrg4_annotated.xlsx"
$ws.Range("H70").Value = "If this is confirmed over benchmark dataset this could be a game changer."
